$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXHome")

# Row 9: Base/Quote feature data
$ws.Range("B9").Value = "Spot"
$ws.Range("C9").Value = "Trader01@Tinyex"
$ws.Range("E9").Value = "USDT"
$ws.Range("D9").Value = "ETH"
$ws.Range("F9:L9").Value = "NA"

$ws.Range("C8").Copy()
$ws.Range("C9").PasteSpecial(-4122)

$ws.Range("F8:L8").Copy()
$ws.Range("F9:L9").PasteSpecial(-4122)

$ws.Range("C8").Copy()
$ws.Range("D9:E9").PasteSpecial(-4122)
$ws.Range("D9:E9").Borders.Item(5).LineStyle = -4142
$ws.Range("D9:E9").Borders.Item(6).LineStyle = -4142
$ws.Range("D9:E9").Borders.Item(7).LineStyle = -4142
$ws.Range("D9:E9").Borders.Item(8).LineStyle = -4142
$ws.Range("D9:E9").Interior.Pattern = -4142

$ws.Range("E11").Select()
